# Update NATMI ligand-receptor pair values for Lrpap1-Ldlr following Dr Hou advice.
# "Ligand-expressing cells" (E) and "Receptor-expressing cells" (K) change from 1 to 3,
# and the dependent expression/specificity statistics are updated to the recomputed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 6.660188000000001
$ws.Range("H2").Value = 19.980564
$ws.Range("I2").Value = 0.1500148400131262
$ws.Range("J2").Value = 0.1500148400131261
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.398836333333334
$ws.Range("N2").Value = 7.196509000000001
$ws.Range("O2").Value = 0.09386760623633866
$ws.Range("P2").Value = 0.09386760623633865
$ws.Range("Q2").Value = 15.97670096123067
$ws.Range("R2").Value = 143.790308651076
$ws.Range("S2").Value = 0.01408153393195947
$ws.Range("T2").Value = 0.01408153393195946

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 6.660188000000001
$ws.Range("H3").Value = 19.980564
$ws.Range("I3").Value = 0.1500148400131262
$ws.Range("J3").Value = 0.1500148400131261
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.682092999999999
$ws.Range("N3").Value = 17.046279
$ws.Range("O3").Value = 0.2223430006085962
$ws.Range("P3").Value = 0.2223430006085962
$ws.Range("Q3").Value = 37.843807613484
$ws.Range("R3").Value = 340.594268521356
$ws.Range("S3").Value = 0.03335474966433698
$ws.Range("T3").Value = 0.03335474966433698

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 6.660188000000001
$ws.Range("H4").Value = 19.980564
$ws.Range("I4").Value = 0.1500148400131262
$ws.Range("J4").Value = 0.1500148400131261
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 12.17451333333334
$ws.Range("N4").Value = 36.52354
$ws.Range("O4").Value = 0.4763944950360188
$ws.Range("P4").Value = 0.4763944950360188
$ws.Range("Q4").Value = 81.08454760850668
$ws.Range("R4").Value = 729.7609284765601
$ws.Range("S4").Value = 0.07146624395596239
$ws.Range("T4").Value = 0.07146624395596238

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 6.660188000000001
$ws.Range("H5").Value = 19.980564
$ws.Range("I5").Value = 0.1500148400131262
$ws.Range("J5").Value = 0.1500148400131261
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.300086333333334
$ws.Range("N5").Value = 15.900259
$ws.Range("O5").Value = 0.2073948981190463
$ws.Range("P5").Value = 0.2073948981190463
$ws.Range("Q5").Value = 35.29957139623068
$ws.Range("R5").Value = 317.696142566076
$ws.Range("S5").Value = 0.03111231246086734
$ws.Range("T5").Value = 0.03111231246086733

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.410331
$ws.Range("H6").Value = 49.230993
$ws.Range("I6").Value = 0.3696281815959916
$ws.Range("J6").Value = 0.3696281815959916
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.398836333333334
$ws.Range("N6").Value = 7.196509000000001
$ws.Range("O6").Value = 0.09386760623633866
$ws.Range("P6").Value = 0.09386760623633865
$ws.Range("Q6").Value = 39.36569824482633
$ws.Range("R6").Value = 354.291284203437
$ws.Range("S6").Value = 0.03469611260390643
$ws.Range("T6").Value = 0.03469611260390642

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.410331
$ws.Range("H7").Value = 49.230993
$ws.Range("I7").Value = 0.3696281815959916
$ws.Range("J7").Value = 0.3696281815959916
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.682092999999999
$ws.Range("N7").Value = 17.046279
$ws.Range("O7").Value = 0.2223430006085962
$ws.Range("P7").Value = 0.2223430006085962
$ws.Range("Q7").Value = 93.24502690278298
$ws.Range("R7").Value = 839.2052421250469
$ws.Range("S7").Value = 0.0821842390055519
$ws.Range("T7").Value = 0.0821842390055519

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 16.410331
$ws.Range("H8").Value = 49.230993
$ws.Range("I8").Value = 0.3696281815959916
$ws.Range("J8").Value = 0.3696281815959916
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 12.17451333333334
$ws.Range("N8").Value = 36.52354
$ws.Range("O8").Value = 0.4763944950360188
$ws.Range("P8").Value = 0.4763944950360188
$ws.Range("Q8").Value = 199.7877935639134
$ws.Range("R8").Value = 1798.09014207522
$ws.Range("S8").Value = 0.1760888309225043
$ws.Range("T8").Value = 0.1760888309225043

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 16.410331
$ws.Range("H9").Value = 49.230993
$ws.Range("I9").Value = 0.3696281815959916
$ws.Range("J9").Value = 0.3696281815959916
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 5.300086333333334
$ws.Range("N9").Value = 15.900259
$ws.Range("O9").Value = 0.2073948981190463
$ws.Range("P9").Value = 0.2073948981190463
$ws.Range("Q9").Value = 86.97617105857634
$ws.Range("R9").Value = 782.785539527187
$ws.Range("S9").Value = 0.07665899906402904
$ws.Range("T9").Value = 0.07665899906402902

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.478895
$ws.Range("H10").Value = 43.436685
$ws.Range("I10").Value = 0.3261242951387937
$ws.Range("J10").Value = 0.3261242951387937
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.398836333333334
$ws.Range("N10").Value = 7.196509000000001
$ws.Range("O10").Value = 0.09386760623633866
$ws.Range("P10").Value = 0.09386760623633865
$ws.Range("Q10").Value = 34.73249939251833
$ws.Range("R10").Value = 312.592494532665
$ws.Range("S10").Value = 0.03061250692019178
$ws.Range("T10").Value = 0.03061250692019177

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 14.478895
$ws.Range("H11").Value = 43.436685
$ws.Range("I11").Value = 0.3261242951387937
$ws.Range("J11").Value = 0.3261242951387937
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.682092999999999
$ws.Range("N11").Value = 17.046279
$ws.Range("O11").Value = 0.2223430006085962
$ws.Range("P11").Value = 0.2223430006085962
$ws.Range("Q11").Value = 82.27042792723499
$ws.Range("R11").Value = 740.4338513451149
$ws.Range("S11").Value = 0.07251145435252283
$ws.Range("T11").Value = 0.07251145435252282

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 14.478895
$ws.Range("H12").Value = 43.436685
$ws.Range("I12").Value = 0.3261242951387937
$ws.Range("J12").Value = 0.3261242951387937
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 12.17451333333334
$ws.Range("N12").Value = 36.52354
$ws.Range("O12").Value = 0.4763944950360188
$ws.Range("P12").Value = 0.4763944950360188
$ws.Range("Q12").Value = 176.2735002294334
$ws.Range("R12").Value = 1586.4615020649
$ws.Range("S12").Value = 0.1553638189016232
$ws.Range("T12").Value = 0.1553638189016232

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 14.478895
$ws.Range("H13").Value = 43.436685
$ws.Range("I13").Value = 0.3261242951387937
$ws.Range("J13").Value = 0.3261242951387937
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 5.300086333333334
$ws.Range("N13").Value = 15.900259
$ws.Range("O13").Value = 0.2073948981190463
$ws.Range("P13").Value = 0.2073948981190463
$ws.Range("Q13").Value = 76.73939351126833
$ws.Range("R13").Value = 690.654541601415
$ws.Range("S13").Value = 0.06763651496445591
$ws.Range("T13").Value = 0.0676365149644559

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 6.847447
$ws.Range("H14").Value = 20.542341
$ws.Range("I14").Value = 0.1542326832520885
$ws.Range("J14").Value = 0.1542326832520885
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 2.398836333333334
$ws.Range("N14").Value = 7.196509000000001
$ws.Range("O14").Value = 0.09386760623633866
$ws.Range("P14").Value = 0.09386760623633865
$ws.Range("Q14").Value = 16.42590465417434
$ws.Range("R14").Value = 147.833141887569
$ws.Range("S14").Value = 0.01447745278028099
$ws.Range("T14").Value = 0.01447745278028098

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 6.847447
$ws.Range("H15").Value = 20.542341
$ws.Range("I15").Value = 0.1542326832520885
$ws.Range("J15").Value = 0.1542326832520885
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.682092999999999
$ws.Range("N15").Value = 17.046279
$ws.Range("O15").Value = 0.2223430006085962
$ws.Range("P15").Value = 0.2223430006085962
$ws.Range("Q15").Value = 38.90783066657099
$ws.Range("R15").Value = 350.170475999139
$ws.Range("S15").Value = 0.03429255758618454
$ws.Range("T15").Value = 0.03429255758618454

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 6.847447
$ws.Range("H16").Value = 20.542341
$ws.Range("I16").Value = 0.1542326832520885
$ws.Range("J16").Value = 0.1542326832520885
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 12.17451333333334
$ws.Range("N16").Value = 36.52354
$ws.Range("O16").Value = 0.4763944950360188
$ws.Range("P16").Value = 0.4763944950360188
$ws.Range("Q16").Value = 83.36433480079334
$ws.Range("R16").Value = 750.2790132071401
$ws.Range("S16").Value = 0.07347560125592892
$ws.Range("T16").Value = 0.07347560125592892

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 6.847447
$ws.Range("H17").Value = 20.542341
$ws.Range("I17").Value = 0.1542326832520885
$ws.Range("J17").Value = 0.1542326832520885
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 5.300086333333334
$ws.Range("N17").Value = 15.900259
$ws.Range("O17").Value = 0.2073948981190463
$ws.Range("P17").Value = 0.2073948981190463
$ws.Range("Q17").Value = 36.29206026292434
$ws.Range("R17").Value = 326.628542366319
$ws.Range("S17").Value = 0.03198707162969403
$ws.Range("T17").Value = 0.03198707162969402
